$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "26/06/2025"
$ws.Range("C2").Value = "114.77 kg (8 cx)"
$ws.Range("D2").Value = "118.97 kg"
$ws.Range("E2").Value = "133.59 kg"
$ws.Range("F2").Value = "2 dias"
